$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.081.42"
$ws.Range("E2").Value = "  -1.92%  "
$ws.Range("D3").Value = "1.835.52"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "323.74"
$ws.Range("E5").Value = "  -3.17%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "0.4634"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").Value = "0.3882"
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("D9").Value = "0.07869"
$ws.Range("E9").Value = "  -0.49%  "
$ws.Range("D10").Value = "0.9638"
$ws.Range("E10").Value = "  -2.07%  "
$ws.Range("D11").Value = "'22.00"
$ws.Range("E11").Value = "  -1.31%  "
$ws.Range("D12").Value = "1.813.05"
$ws.Range("E12").Value = "  -1.32%  "
$ws.Range("E13").Value = "  -2.56%  "
$ws.Range("D14").Value = "6.927"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").Value = "0.06842"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "'88.50"
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "0.000009958"
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("D19").Value = "16.74"
$ws.Range("E19").Value = "  -2.27%  "
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "28.107.37"
$ws.Range("E21").Value = "  -1.88%  "
$ws.Range("D22").Value = "5.323"
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("D23").Value = "11.02"
$ws.Range("E23").Value = "  -2.58%  "
$ws.Range("D24").Value = "2.099"
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("D25").Value = "2.041.10"
$ws.Range("E25").Value = "  -2.05%  "
$ws.Range("D26").Value = "154.77"
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("D27").Value = "'19.20"
$ws.Range("E27").Value = "  -1.47%  "
$ws.Range("D28").Value = "5.701"
$ws.Range("E28").Value = "  -5.64%  "
$ws.Range("D29").Value = "1.969"
$ws.Range("E29").Value = "  -3.04%  "
$ws.Range("D30").Value = "118.28"
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("D31").Value = "0.9364"
$ws.Range("E31").Value = "  -4.30%  "
$ws.Range("D32").Value = "0.09234"
$ws.Range("E32").Value = "  -1.74%  "
$ws.Range("D33").Value = "5.276"
$ws.Range("E33").Value = "  -1.83%  "
$ws.Range("E34").Value = "  -1.91%  "
$ws.Range("D35").Value = "3.307"
$ws.Range("E35").Value = "  -5.05%  "
$ws.Range("D36").Value = "0.05879"
$ws.Range("E36").Value = "  -4.45%  "
$ws.Range("D37").Value = "0.02131"
$ws.Range("E37").Value = "  -3.09%  "
$ws.Range("D38").Value = "1.146"
$ws.Range("E38").Value = "  -2.02%  "
$ws.Range("D39").Value = "'7.790"
$ws.Range("E39").Value = "  +2.35%  "
$ws.Range("D40").Value = "'0.5610"
$ws.Range("E40").Value = "  -1.92%  "
$ws.Range("D41").Value = "9.925"
$ws.Range("E41").Value = "  -2.65%  "
$ws.Range("D42").Value = "0.1766"
$ws.Range("E42").Value = "  -1.96%  "
$ws.Range("E43").Value = "  +1.68%  "
$ws.Range("E44").Value = "  -0.62%  "
$ws.Range("D45").Value = "0.5285"
$ws.Range("E45").Value = "  -1.96%  "
$ws.Range("D46").Value = "1.159"
$ws.Range("E46").Value = "  -6.88%  "
$ws.Range("D47").Value = "2.122"
$ws.Range("E47").Value = "  -10.63%  "
$ws.Range("D48").Value = "1.827"
$ws.Range("E48").Value = "  -4.21%  "
$ws.Range("D49").Value = "112.52"
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").Value = "1.001"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "1.028"
$ws.Range("E51").Value = "  +0.77%  "
